# Applies the "Anzahl Kinder" table + employee sheet changes described in the diff.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Tabelle1 (Max Mustermann data)
$ws2 = $wb.Worksheets.Item(2)   # Tabelle2 (lookup lists)

# ---------------------------------------------------------------------------
# 1. Tabelle1 ("Max Mustermann") - rewrite rows 48-56 so that a new
#    "Anzahl Kinder" entry is inserted right after "Zusatzbeitrag
#    Krankenversicherung in Prozent" (row 50) and the rest of the rows that
#    used to belong to the "wohnhaft Sachsen?" block shift down by one.
# ---------------------------------------------------------------------------

# D48:D54 switch their highlight color from style 11 (red) to style 10 (green)
# -> copy the format that D2 already uses (style 10) onto D48:D54.
$ws1.Range("D2").Copy() | Out-Null
$ws1.Range("D48:D54").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 50 no longer is the end of its little section (that moved to row 56),
# so it loses the custom row height / thick bottom border.
$ws1.Rows.Item(50).AutoFit() | Out-Null

# New content for A51:B51 ("Anzahl Kinder") - give it the plain style (8)
# that the surrounding rows use (copy format from A53 / B53 which already
# carry style 8) before writing the values.
$ws1.Range("A53").Copy() | Out-Null
$ws1.Range("A51").PasteSpecial(-4122) | Out-Null
$ws1.Range("B53").Copy() | Out-Null
$ws1.Range("B51").PasteSpecial(-4122) | Out-Null

$ws1.Range("A51").Value = "Anzahl Kinder"
$ws1.Range("B51").Value = 2

# Row 52 -> "AN-Pflegeversicherungsbeitrag in Prozent" = 1
$ws1.Range("A52").Value = "AN-Pflegeversicherungsbeitrag in Prozent"
$ws1.Range("B52").Value = 1

# Row 53 -> "Beitragsbemessungsgrenze Pflegeversicherung Ost" = 35000
$ws1.Range("A53").Value = "Beitragsbemessungsgrenze Pflegeversicherung Ost"
$ws1.Range("B53").Value = 35000

# Row 54 -> "Beitragsbemessungsgrenze Pflegeversicherung West" = 38000
$ws1.Range("A54").Value = "Beitragsbemessungsgrenze Pflegeversicherung West"
$ws1.Range("B54").Value = 38000

# Row 55 -> "wohnhaft Sachsen?" (keeps the plain style it already has)
$ws1.Range("A55").Value = "wohnhaft Sachsen?"

# Row 56 -> "AG-Pflegeversicherungsbeitrag in Prozent" (keeps its own
# bottom-border style which marks the end of the section)
$ws1.Range("A56").Value = "AG-Pflegeversicherungsbeitrag in Prozent"

# ---------------------------------------------------------------------------
# 2. Data validation on Tabelle1: B51 used to be part of the "ja/nein" list
#    (together with B52); now B51 gets its own validation against the new
#    "Anzahl Kinder" list on Tabelle2 (I2:I102), while B52 has no validation
#    list at all anymore (it is a free numeric input).
# ---------------------------------------------------------------------------
$ws1.Range("B51:B52").Validation.Delete()
$ws1.Range("B51").Validation.Add(3, 1, 1, "=Tabelle2!`$I`$2:`$I`$102") | Out-Null

# ---------------------------------------------------------------------------
# 3. Tabelle1 sheet view - scroll position / selection changed as part of
#    the edit session.
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 31
$ws1.Range("B55").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Tabelle2 - add the new "Anzahl Kinder" lookup column (I1:I102) with
#    values 0..100.
# ---------------------------------------------------------------------------
$ws2.Range("I1").Value = "Anzahl Kinder"
$ws2.Range("I1").Font.Bold = $true

$values = New-Object 'object[,]' 101,1
for ($i = 0; $i -le 100; $i++) {
    $values[$i,0] = $i
}
$ws2.Range("I2:I102").Value = $values

# ---------------------------------------------------------------------------
# 5. Turn the five lookup ranges on Tabelle2 into real Excel Tables
#    (ListObjects), matching the tables added in the diff.
# ---------------------------------------------------------------------------
$tbl1 = $ws2.ListObjects.Add(1, $ws2.Range("A1:A4"), [System.Type]::Missing, 1)
$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("C1:C7"), [System.Type]::Missing, 1)
$tbl3 = $ws2.ListObjects.Add(1, $ws2.Range("E1:E7"), [System.Type]::Missing, 1)
$tbl4 = $ws2.ListObjects.Add(1, $ws2.Range("G1:G3"), [System.Type]::Missing, 1)
$tbl5 = $ws2.ListObjects.Add(1, $ws2.Range("I1:I102"), [System.Type]::Missing, 1)

$ws2.ListObjects.Item(5).Name = "tbl_Anzahl_Kinder"
$ws2.ListObjects.Item(4).Name = "tbl_ja_nein"
$ws2.ListObjects.Item(3).Name = "tbl_Steuerklasse"
$ws2.ListObjects.Item(2).Name = "tbl_Mitarbeitertyp"
$ws2.ListObjects.Item(1).Name = "tbl_Geschlecht"

$ws2.ListObjects.Item(5).TableStyle = "TableStyleMedium9"
$ws2.ListObjects.Item(4).TableStyle = "TableStyleMedium9"
$ws2.ListObjects.Item(3).TableStyle = "TableStyleMedium9"
$ws2.ListObjects.Item(2).TableStyle = "TableStyleMedium9"
$ws2.ListObjects.Item(1).TableStyle = "TableStyleMedium9"

# ---------------------------------------------------------------------------
# 6. Tabelle2 column widths / selection.
# ---------------------------------------------------------------------------
$ws2.Range("A1").EntireColumn.ColumnWidth = 11.666666666666666
$ws2.Range("C1").EntireColumn.ColumnWidth = 14.833333333333332
$ws2.Range("E1").EntireColumn.ColumnWidth = 13.0
$ws2.Range("I1").EntireColumn.ColumnWidth = 14.0

$ws2.Activate()
$ws2.Range("I1").Select() | Out-Null

$ws1.Activate()
